# Pulse_Manager.xlsx edit:
#  - Pulse_Queue: rename "End Time" header to "Pulse Length", add a new
#    queue row for Pulse_5 (320 -> 400).
#  - Pulse_Definition: add a new "Pulse_5" column (F) with its channel
#    settings, fix a few existing On/Off values, extend the data
#    validation + column formatting to the new column, and make
#    Pulse_Definition the active/visible sheet.

$wb = $excel.ActiveWorkbook

$wsQueue = $wb.Worksheets.Item("Pulse_Queue")
$wsDef   = $wb.Worksheets.Item("Pulse_Definition")

# ---------------------------------------------------------------------
# Pulse_Queue: relabel the third column and append the Pulse_5 entry
# ---------------------------------------------------------------------
$wsQueue.Range("C1").Value = "Pulse Length"

$wsQueue.Range("A7").Value = "Pulse_5"
$wsQueue.Range("B7").Value = 320
$wsQueue.Range("C7").Value = 400

# ---------------------------------------------------------------------
# Pulse_Definition: correct a handful of existing channel settings
# ---------------------------------------------------------------------
$wsDef.Range("B3").Value = "On"      # Ch1 / Pulse_1 : Sweep -> On
$wsDef.Range("C4").Value = "On"      # Ch2 / Pulse_2 : Off   -> On
$wsDef.Range("B6").Value = "On"      # Ch4 / Pulse_1 : Off   -> On

# ---------------------------------------------------------------------
# Pulse_Definition: new Pulse_5 column (F)
# ---------------------------------------------------------------------
$wsDef.Range("F1").Value = "Pulse_5"
$wsDef.Range("F2").Value = "Off"
$wsDef.Range("F3").Value = "Sweep"
$wsDef.Range("F4").Value = "Off"
$wsDef.Range("F5").Value = "On"
$wsDef.Range("F6").Value = "Off"
$wsDef.Range("F7:F29").Value = "Off"

# match the existing data columns' width for the new column
$wsDef.Columns.Item(6).ColumnWidth = 9.3

# data validation (list On/Sweep/Off) needs to cover the new column too
$wsDef.Range("B2:E29").Validation.Delete() | Out-Null
$wsDef.Range("B2:F29").Validation.Add(3, 1, 1, '"On, Sweep, Off"') | Out-Null

# ---------------------------------------------------------------------
# Selection / active sheet bookkeeping
# ---------------------------------------------------------------------
$wsQueue.Activate() | Out-Null
$wsQueue.Range("C8").Select() | Out-Null

$wsDef.Activate() | Out-Null
$wsDef.Range("G8").Select() | Out-Null
